$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 907
$ws.Range("F3").Value = 1021
$ws.Range("F4").Value = 804
$ws.Range("F5").Value = 879
$ws.Range("F7").Value = 701
$ws.Range("F9").Value = 1307
$ws.Range("F10").Value = 726
$ws.Range("F12").Value = 555
$ws.Range("F13").Value = 188
$ws.Range("F14").Value = 49
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 1099
$ws.Range("F17").Value = 88
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 420
$ws.Range("F22").Value = 597
$ws.Range("F23").Value = 154
$ws.Range("F24").Value = 654
$ws.Range("F25").Value = 37
$ws.Range("F26").Value = 1062

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 344
$ws.Range("F3").Value = 113
$ws.Range("F6").Value = 191
$ws.Range("F7").Value = 249
$ws.Range("F11").Value = 114

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 380

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 380
$ws.Range("F3").Value = 344
$ws.Range("F4").Value = 907
$ws.Range("F5").Value = 1021
$ws.Range("F6").Value = 804
$ws.Range("F7").Value = 879
$ws.Range("F9").Value = 701
$ws.Range("F11").Value = 1307
$ws.Range("F12").Value = 726
$ws.Range("F13").Value = 113
$ws.Range("F16").Value = 555
$ws.Range("F18").Value = 188
$ws.Range("F19").Value = 49
$ws.Range("F20").Value = 49
$ws.Range("F21").Value = 1099
$ws.Range("F22").Value = 191
$ws.Range("F23").Value = 89
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 420
$ws.Range("F28").Value = 249
$ws.Range("F30").Value = 597
$ws.Range("F33").Value = 114
$ws.Range("F34").Value = 114
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 654
$ws.Range("F37").Value = 37
$ws.Range("F38").Value = 1062
